# CourseInformation.xlsx - add remaining lower-division Physics courses
# required for CS majors, and backfill the Difficulty column for the
# courses that already existed (E8:E14), matching wrap-text formatting
# that was extended to column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill "Difficulty" (column E) for the pre-existing rows 8-14 ---
$ws.Range("E8").Value = 3
$ws.Range("E9").Value = 3
$ws.Range("E10").Value = 4
$ws.Range("E11").Value = 4
$ws.Range("E12").Value = 4
$ws.Range("E13").Value = 5
$ws.Range("E14").Value = 4

# --- New course rows (Physics lower-division requirements) ---
# Column A (course names) written first, in row order, so new shared
# strings are interned in the same order the workbook author typed them.
$ws.Range("A15").Value = "Physics 1A"
$ws.Range("A16").Value = "Physics 1B"
$ws.Range("A17").Value = "Physics 1C"
$ws.Range("A18").Value = "Physics 4AL"
$ws.Range("A19").Value = "Physics 4BL"

# Column B / C for all five new rows.
$ws.Range("B15").Value = "Lower Division Requirement"
$ws.Range("C15").Value = 5
$ws.Range("B16").Value = "Lower Division Requirement"
$ws.Range("C16").Value = 5
$ws.Range("B17").Value = "Lower Division Requirement"
$ws.Range("C17").Value = 5
$ws.Range("B18").Value = "Lower Division Requirement"
$ws.Range("C18").Value = 2
$ws.Range("B19").Value = "Lower Division Requirement"
$ws.Range("C19").Value = 2

# Column D (descriptions) for rows 15-17.
$ws.Range("D15").Value = "Motion, Newton laws, work, energy, linear and angular momentum, rotation, equilibrium, gravitation."
$ws.Range("D16").Value = "Fluid mechanics, oscillation, mechanical waves, and sound. Electric charge, field and potential, capacitors, and dielectrics. Currents and resistance, direct-current circuits."
$ws.Range("D17").Value = "Magnetic fields, Ampere's law, Faraday's law, inductance, and alternating current circuits. Maxwell's equations, electromagnetic waves, light, geometrical optics, interference and diffraction."

# Column F (prereqs) for row 15 (reuses an existing shared string) and
# rows 16-17 (new strings).
$ws.Range("F15").Value = "Math 31A, Math 31B, Math 32A"
$ws.Range("F16").Value = "Math 31A, Math 31B, Math 32A, Math 32B, Physics 1A"
$ws.Range("F17").Value = "Math 31A, Math 31B, Math 32A, Math 32B, Math 33A, Physics 1A, Physics 1B"

# Column D (descriptions) for rows 18-19.
$ws.Range("D18").Value = "*Fulfills Same Requirement as Physics 4BL* Computerized measurements of uniform and accelerated motion, including oscillations. Analysis of data and comparison of results to predictions, including least-squares fitting. Conception, execution, and presentation of creative projects involving motion."
$ws.Range("D19").Value = "*Fulfills Same Requirement as Physics 4AL* Sound waves and electric circuits, taken by digital oscilloscopes and analyzed by Fourier transformation. Geometrical and physical optics. Conception, execution, and presentation of creative projects involving sound waves or electric circuits."

# Column F (prereqs) for rows 18-19.
$ws.Range("F18").Value = "Math 31A, Math 31B, Math 32A, Physics 1A, Physics 1B"
$ws.Range("F19").Value = "Math 31A, Math 31B, Math 32A, Math 32B, Math 33A, Physics 1A, Physics 1B, Physics 1C"

# --- Formatting ---
# Wrap text on the whole Prereqs column (F), matching the new column
# style and the per-cell style applied to F1:F19.
$ws.Range("F1:F19").WrapText = $true

# Row heights for the new rows (wrapped text autosizing).
$ws.Rows.Item(15).RowHeight = 51
$ws.Rows.Item(16).RowHeight = 68
$ws.Rows.Item(17).RowHeight = 68
$ws.Rows.Item(18).RowHeight = 104
$ws.Rows.Item(19).RowHeight = 102

# --- View state ---
$ws.Range("F20").Select()
